$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04763786555579896
$ws.Range("C2").Value = 0.04240448674262143
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 9.555764474831394

$ws.Range("B3").Value = 0.6753301551942219
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 9844.520545567508
$ws.Range("E3").Value = 616238.5361209477
$ws.Range("G3").Value = 626084.0447870663
